$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 8599.625
$ws.Cells.Item(18, 9).Value = 8399.571
$ws.Cells.Item(18, 11).Value = 8399.571
$ws.Cells.Item(18, 13).Value = -8115.571

$ws.Cells.Item(32, 8).Value = 2877.0833
$ws.Cells.Item(32, 9).Value = 1287.6666
$ws.Cells.Item(32, 10).Value = 4466.5
$ws.Cells.Item(32, 11).Value = 1287.6666
$ws.Cells.Item(32, 12).Value = 4466.5
$ws.Cells.Item(32, 13).Value = -961.6666
$ws.Cells.Item(32, 14).Value = -5118.5

$ws.Cells.Item(53, 8).Value = 427.6
$ws.Cells.Item(53, 9).Value = 280.5
$ws.Cells.Item(53, 10).Value = 525.6667
$ws.Cells.Item(53, 11).Value = 280.5
$ws.Cells.Item(53, 12).Value = 525.6667
$ws.Cells.Item(53, 13).Value = 356.5
$ws.Cells.Item(53, 14).Value = -1799.6667

$ws.Cells.Item(111, 8).Value = 795.5
$ws.Cells.Item(111, 10).Value = 795
$ws.Cells.Item(111, 12).Value = 2385
$ws.Cells.Item(111, 14).Value = -8519

$ws.Cells.Item(116, 8).Value = 4899.4
$ws.Cells.Item(116, 10).Value = 3499.3333
$ws.Cells.Item(116, 12).Value = 3499.3333
$ws.Cells.Item(116, 14).Value = -10383.3333

$ws.Cells.Item(141, 8).Value = 2873.4736
$ws.Cells.Item(141, 9).Value = 2820.6667
$ws.Cells.Item(141, 10).Value = 3071.5
$ws.Cells.Item(141, 11).Value = 8462.000100000001
$ws.Cells.Item(141, 12).Value = 9214.5
$ws.Cells.Item(141, 13).Value = -3282.000100000001
$ws.Cells.Item(141, 14).Value = -19574.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14, 8).Value = 2425
$ws.Cells.Item(14, 9).Value = 2425
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 2425
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -2250
$ws.Cells.Item(14, 14).Value = ""

$ws.Cells.Item(61, 8).Value = 3099.8333
$ws.Cells.Item(61, 9).Value = 2999.8
$ws.Cells.Item(61, 11).Value = 2999.8
$ws.Cells.Item(61, 13).Value = -2787.8

$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).Value = ""

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).Value = ""

$ws.Cells.Item(136, 8).Value = 3099.8333
$ws.Cells.Item(136, 9).Value = 2999.8
$ws.Cells.Item(136, 11).Value = 8999.400000000001
$ws.Cells.Item(136, 13).Value = -6449.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1911.4
$ws.Cells.Item(20, 9).Value = 1139.25
$ws.Cells.Item(20, 10).Value = 5000
$ws.Cells.Item(20, 11).Value = 1139.25
$ws.Cells.Item(20, 12).Value = 5000
$ws.Cells.Item(20, 13).Value = -892.25
$ws.Cells.Item(20, 14).Value = -5494

$ws.Cells.Item(105, 8).Value = 1759.5
$ws.Cells.Item(105, 9).Value = 1788.6923
$ws.Cells.Item(105, 10).Value = 1633
$ws.Cells.Item(105, 11).Value = 1788.6923
$ws.Cells.Item(105, 12).Value = 1633
$ws.Cells.Item(105, 13).Value = -41.69229999999993
$ws.Cells.Item(105, 14).Value = -5127

$ws.Cells.Item(107, 8).Value = 1513.44
$ws.Cells.Item(107, 9).Value = 1749.579
$ws.Cells.Item(107, 10).Value = 765.6667
$ws.Cells.Item(107, 11).Value = 1749.579
$ws.Cells.Item(107, 12).Value = 765.6667
$ws.Cells.Item(107, 13).Value = 170.421
$ws.Cells.Item(107, 14).Value = -4605.6667

$ws.Cells.Item(134, 8).Value = 2512.1667
$ws.Cells.Item(134, 9).Value = 2512.1667
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 7536.500100000001
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -5001.500100000001
$ws.Cells.Item(134, 14).Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 6499.75
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 6499.75
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 6499.75
$ws.Cells.Item(12, 13).Value = ""
$ws.Cells.Item(12, 14).Value = -6839.75

$ws.Cells.Item(16, 8).Value = 7247.5
$ws.Cells.Item(16, 9).Value = 6330
$ws.Cells.Item(16, 11).Value = 6330
$ws.Cells.Item(16, 13).Value = -6043

$ws.Cells.Item(31, 8).Value = 1523
$ws.Cells.Item(31, 9).Value = 1030.6666
$ws.Cells.Item(31, 11).Value = 1030.6666
$ws.Cells.Item(31, 13).Value = -735.6666

$ws.Cells.Item(34, 8).Value = 1523
$ws.Cells.Item(34, 9).Value = 1030.6666
$ws.Cells.Item(34, 11).Value = 1030.6666
$ws.Cells.Item(34, 13).Value = -828.6666

$ws.Cells.Item(113, 8).Value = 7247.5
$ws.Cells.Item(113, 9).Value = 6330
$ws.Cells.Item(113, 11).Value = 6330
$ws.Cells.Item(113, 13).Value = -4160

$ws.Cells.Item(134, 8).Value = 3205.5
$ws.Cells.Item(134, 9).Value = 3028.5715
$ws.Cells.Item(134, 11).Value = 9085.7145
$ws.Cells.Item(134, 13).Value = -6550.7145

$ws.Cells.Item(141, 8).Value = 166994.62
$ws.Cells.Item(141, 10).Value = 166994.62
$ws.Cells.Item(141, 12).Value = 166994.62
$ws.Cells.Item(141, 14).Value = -177354.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 160.6
$ws.Cells.Item(2, 9).Value = 56.5
$ws.Cells.Item(2, 10).Value = 230
$ws.Cells.Item(2, 11).Value = 339
$ws.Cells.Item(2, 12).Value = 1380
$ws.Cells.Item(2, 13).Value = -226
$ws.Cells.Item(2, 14).Value = -1606

$ws.Cells.Item(33, 8).Value = 1774.25
$ws.Cells.Item(33, 9).Value = 1774.25
$ws.Cells.Item(33, 11).Value = 10645.5
$ws.Cells.Item(33, 13).Value = -10362.5

$ws.Cells.Item(40, 8).Value = 60.125
$ws.Cells.Item(40, 9).Value = 47
$ws.Cells.Item(40, 10).Value = 99.5
$ws.Cells.Item(40, 11).Value = 188
$ws.Cells.Item(40, 12).Value = 398
$ws.Cells.Item(40, 13).Value = -119
$ws.Cells.Item(40, 14).Value = -536

$ws.Cells.Item(56, 8).Value = 11173.3955
$ws.Cells.Item(56, 9).Value = 11173.3955
$ws.Cells.Item(56, 11).Value = 11173.3955
$ws.Cells.Item(56, 13).Value = -10643.3955

$ws.Cells.Item(81, 8).Value = 2266.6667
$ws.Cells.Item(81, 9).Value = 2000
$ws.Cells.Item(81, 10).Value = 2400
$ws.Cells.Item(81, 11).Value = 6000
$ws.Cells.Item(81, 12).Value = 7200
$ws.Cells.Item(81, 13).Value = -4877
$ws.Cells.Item(81, 14).Value = -9446

$ws.Cells.Item(84, 8).Value = 2266.6667
$ws.Cells.Item(84, 9).Value = 2000
$ws.Cells.Item(84, 10).Value = 2400
$ws.Cells.Item(84, 11).Value = 18000
$ws.Cells.Item(84, 12).Value = 21600
$ws.Cells.Item(84, 13).Value = -12384
$ws.Cells.Item(84, 14).Value = -32832

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6366.6665
$ws.Cells.Item(70, 9).Value = 6350
$ws.Cells.Item(70, 11).Value = 6350
$ws.Cells.Item(70, 13).Value = -6080

$ws.Cells.Item(73, 8).Value = 6366.6665
$ws.Cells.Item(73, 9).Value = 6350
$ws.Cells.Item(73, 11).Value = 6350
$ws.Cells.Item(73, 13).Value = -5414

$ws.Cells.Item(107, 8).Value = 2378.6365
$ws.Cells.Item(107, 9).Value = 1069.6666
$ws.Cells.Item(107, 10).Value = 2869.5
$ws.Cells.Item(107, 11).Value = 1069.6666
$ws.Cells.Item(107, 12).Value = 2869.5
$ws.Cells.Item(107, 13).Value = 850.3334
$ws.Cells.Item(107, 14).Value = -6709.5

$ws.Cells.Item(113, 8).Value = 249
$ws.Cells.Item(113, 9).Value = 249
$ws.Cells.Item(113, 11).Value = 249
$ws.Cells.Item(113, 13).Value = 1921

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2897.6667
$ws.Cells.Item(40, 9).Value = 2337.2
$ws.Cells.Item(40, 11).Value = 2337.2
$ws.Cells.Item(40, 13).Value = -2201.2

$ws.Cells.Item(46, 8).Value = 3292.5
$ws.Cells.Item(46, 9).Value = 2791
$ws.Cells.Item(46, 10).Value = 5800
$ws.Cells.Item(46, 11).Value = 2791
$ws.Cells.Item(46, 12).Value = 5800
$ws.Cells.Item(46, 13).Value = -2603
$ws.Cells.Item(46, 14).Value = -6176

$ws.Cells.Item(122, 8).Value = 6175.923
$ws.Cells.Item(122, 9).Value = 4672.846
$ws.Cells.Item(122, 10).Value = 7679
$ws.Cells.Item(122, 11).Value = 14018.538
$ws.Cells.Item(122, 12).Value = 23037
$ws.Cells.Item(122, 13).Value = -11568.538
$ws.Cells.Item(122, 14).Value = -27937

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 1506
$ws.Cells.Item(23, 9).Value = 1432.5
$ws.Cells.Item(23, 11).Value = 1432.5
$ws.Cells.Item(23, 13).Value = -1203.5

$ws.Cells.Item(107, 8).Value = 789.75
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 789.75
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 2369.25
$ws.Cells.Item(107, 13).Value = ""
$ws.Cells.Item(107, 14).Value = -6209.25
